# Generate Report for Handoff
# Updates the "Ready for handoff" status (previously "Handed back: in sync
# with en-US") and refreshed timestamps for the 6c2f9e99-*.md /
# b870dd9b-*.md files, across the Overview, zh-cn and de-de sheets, and
# records a "version not latest" error detail for both of those files on
# the zh-cn / de-de sheets. Also widens the Error Detail column.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"
$overviewTimestamp = "2016-08-30 16:31:44"
$zhcnTimestamp = "2016-08-30 16:31:39"
$dedeTimestamp = "2016-08-30 16:31:44"

$msg6c2f9e99 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f565a285417940a0d6151576817cb7852336ad35/e2e/6c2f9e99-5d49-487d-b1b7-c6ffac62813c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f1e5032e4c234021f83b6dce5edb35112d98e9b/e2e/6c2f9e99-5d49-487d-b1b7-c6ffac62813c.md."
$msgb870dd9b = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f565a285417940a0d6151576817cb7852336ad35/e2e/b870dd9b-8bc5-4a17-985d-ac85d59d771d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f1e5032e4c234021f83b6dce5edb35112d98e9b/e2e/b870dd9b-8bc5-4a17-985d-ac85d59d771d.md."

# ---------------------------------------------------------------
# Overview sheet: rows 4 (6c2f9e99-...) and 5 (b870dd9b-...)
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E4").Value = $readyStatus
$overview.Range("F4").Value = $readyStatus
$overview.Range("G4").Value = $overviewTimestamp

$overview.Range("E5").Value = $readyStatus
$overview.Range("F5").Value = $readyStatus
$overview.Range("G5").Value = $overviewTimestamp

# ---------------------------------------------------------------
# zh-cn sheet: rows 4 (6c2f9e99-...) and 5 (b870dd9b-...)
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C4").Value = $readyStatus
$zhcn.Range("H4").Value = $zhcnTimestamp
$zhcn.Range("P4").Value = $msg6c2f9e99

$zhcn.Range("C5").Value = $readyStatus
$zhcn.Range("H5").Value = $zhcnTimestamp
$zhcn.Range("P5").Value = $msgb870dd9b

# Widen the Error Detail column (column 16 / P) to fit the new message
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------
# de-de sheet: rows 4 (6c2f9e99-...) and 5 (b870dd9b-...)
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C4").Value = $readyStatus
$dede.Range("H4").Value = $dedeTimestamp
$dede.Range("P4").Value = $msg6c2f9e99

$dede.Range("C5").Value = $readyStatus
$dede.Range("H5").Value = $dedeTimestamp
$dede.Range("P5").Value = $msgb870dd9b

# Widen the Error Detail column (column 16 / P) to fit the new message
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664

Write-Output "Report generated for handoff"
